# Rename the "MassWateR" organization label to "MassBays" in the Meta sheet,
# and move the selection (active cell) on the sheet to F9 — matching the
# authoring app's state when the workbook was re-saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "MassBays"
$ws.Range("B5").Value = "MassBays"
$ws.Range("B6").Value = "MassBays"

$ws.Range("F9").Select()
